# Contact Us page data-entry sheet update.
# Adds a new "Status" column (S) with Positive/negative values per row,
# tweaks a few of the sample rows (organization name, zip/phone numbers,
# fax-number cell becomes a text value, and one e-mail becomes a short
# name), and scrolls/selects so M3 is the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Status" header column -------------------------------------------------
$ws.Range("S1").Value = "Status"

# --- Row 2 (Contractor / Softway23) ---------------------------------------------
$ws.Range("S2").Value = "Positive"

# --- Row 3 (Developer / was Softway21) -------------------------------------------
# Keep shared-string insertion order lined up with the target workbook:
# "negative" is introduced here before "softway".
$ws.Range("S3").Value = "negative"
$ws.Range("C3").Value = "softway"
$ws.Range("K3").Value = 123456
$ws.Range("L3").Value = 12345645454

# --- Row 5 (HVAC Dealer) e-mail shortened to "rais" before row 3's "test" -------
$ws.Range("N5").Value = "rais"

# Fax Number cell on row 3 becomes the literal text "test" (was numeric).
$ws.Range("M3").Value = "test"

# --- Row 4 (Homeowner / Softway26) -----------------------------------------------
$ws.Range("S4").Value = "Positive"

# --- Row 5 (HVAC Dealer / Softway25) ---------------------------------------------
$ws.Range("S5").Value = "negative"

# --- View state: scroll so column B is leftmost, select M3 ----------------------
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M3").Select()
